# Update "想去人数" (interested-count) figures in the F column across the
# relevant worksheets, matching the upstream data refresh captured in the
# commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 523
$ws1.Range("F4").Value  = 1519
$ws1.Range("F5").Value  = 152
$ws1.Range("F9").Value  = 737
$ws1.Range("F10").Value = 1046
$ws1.Range("F14").Value = 6393
$ws1.Range("F15").Value = 7
$ws1.Range("F18").Value = 151
$ws1.Range("F20").Value = 15290
$ws1.Range("F21").Value = 1518
$ws1.Range("F23").Value = 141
$ws1.Range("F25").Value = 11041
$ws1.Range("F27").Value = 4317
$ws1.Range("F29").Value = 373
$ws1.Range("F31").Value = 302

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 344

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 523
$ws4.Range("F4").Value  = 1519
$ws4.Range("F5").Value  = 152
$ws4.Range("F7").Value  = 344
$ws4.Range("F10").Value = 737
$ws4.Range("F12").Value = 1046
$ws4.Range("F17").Value = 6393
$ws4.Range("F18").Value = 7
$ws4.Range("F21").Value = 151
$ws4.Range("F23").Value = 15290
$ws4.Range("F24").Value = 1518
$ws4.Range("F26").Value = 141
$ws4.Range("F28").Value = 11041
$ws4.Range("F30").Value = 4317
$ws4.Range("F32").Value = 373
$ws4.Range("F34").Value = 302
